$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A407").Value = "Buying Opportunity"
$ws.Range("B407").Value = "support Zone"
$ws.Range("C407").Value = "long buildup"
$ws.Range("D407").Value = "Short buildup"
$ws.Range("E407").Value = "FII ENTERING"

$ws.Range("A408").Value = "AMJLAND"
$ws.Range("B408").Value = "ACC"
$ws.Range("C408").Value = "GODREJCP"
$ws.Range("D408").Value = "ABB"
$ws.Range("E408").Value = "EXIDEIND"
$ws.Range("F408").Value = 41.5
$ws.Range("G408").Value = 2570.55
$ws.Range("H408").Value = 1331
$ws.Range("I408").Value = 8233.65
$ws.Range("J408").Value = 509.9

$ws.Range("A409").Value = "EXIDEIND"
$ws.Range("B409").Value = "AMNPLST"
$ws.Range("C409").Value = "HDFCAMC"
$ws.Range("D409").Value = "RAMCOCEM"
$ws.Range("F409").Value = 509.9
$ws.Range("G409").Value = 189.25
$ws.Range("H409").Value = 4011.6
$ws.Range("I409").Value = 762.85

$ws.Range("A410").Value = "OBCL"
$ws.Range("B410").Value = "ANDHRAPAP"
$ws.Range("C410").Value = "JUBLFOOD"
$ws.Range("F410").Value = 57.9
$ws.Range("G410").Value = 509.1
$ws.Range("H410").Value = 498

$ws.Range("A411").Value = "RUSHIL"
$ws.Range("B411").Value = "AWL"
$ws.Range("F411").Value = 331.7
$ws.Range("G411").Value = 338.25

$ws.Range("B412").Value = "BECTORFOOD"
$ws.Range("G412").Value = 1202.35

$ws.Range("B413").Value = "CARBORUNIV"
$ws.Range("G413").Value = 1603.55

$ws.Range("B414").Value = "CASTROLIND"
$ws.Range("G414").Value = 187.2

$ws.Range("B415").Value = "CELEBRITY"
$ws.Range("G415").Value = 18.95

$ws.Range("B416").Value = "CLEDUCATE"
$ws.Range("G416").Value = 77.1

$ws.Range("B417").Value = "DATAMATICS"
$ws.Range("G417").Value = 549.05

$ws.Range("B418").Value = "DEEPINDS"
$ws.Range("G418").Value = 278.95

$ws.Range("B419").Value = "DELTACORP"
$ws.Range("G419").Value = 115.85

$ws.Range("B420").Value = "DODLA"
$ws.Range("G420").Value = 879.25

$ws.Range("B421").Value = "EXXARO"
$ws.Range("G421").Value = 98.55

$ws.Range("B422").Value = "FAZE3Q"
$ws.Range("G422").Value = 398.6

$ws.Range("B423").Value = "GENCON"
$ws.Range("G423").Value = 39.75

$ws.Range("B424").Value = "GEPIL"
$ws.Range("G424").Value = 328.05

$ws.Range("B425").Value = "GMBREW"
$ws.Range("G425").Value = 635.1

$ws.Range("B426").Value = "GNA"
$ws.Range("G426").Value = 393

$ws.Range("B427").Value = "IFCI"
$ws.Range("G427").Value = 57.25

$ws.Range("B428").Value = "INDIANHUME"
$ws.Range("G428").Value = 318.4

$ws.Range("B429").Value = "INDIGO"
$ws.Range("G429").Value = 4197.05

$ws.Range("B430").Value = "INDORAMA"
$ws.Range("G430").Value = 41.4

$ws.Range("B431").Value = "IRMENERGY"
$ws.Range("G431").Value = 469.3

$ws.Range("B432").Value = "ITDC"
$ws.Range("G432").Value = 662.75

$ws.Range("B433").Value = "JAYSREETEA"
$ws.Range("G433").Value = 99.35

$ws.Range("B434").Value = "KAKATCEM"
$ws.Range("G434").Value = 207.55

$ws.Range("B435").Value = "KMSUGAR"
$ws.Range("G435").Value = 35.9

$ws.Range("B436").Value = "KOTHARIPET"
$ws.Range("G436").Value = 130.95

$ws.Range("B437").Value = "KSCL"
$ws.Range("G437").Value = 850.35

$ws.Range("B438").Value = "LXCHEM"
$ws.Range("G438").Value = 244.85

$ws.Range("B439").Value = "MANGALAM"
$ws.Range("G439").Value = 103.95

$ws.Range("B440").Value = "MANGCHEFER"
$ws.Range("G440").Value = 105.55

$ws.Range("B441").Value = "MGEL"
$ws.Range("G441").Value = 21

$ws.Range("B442").Value = "MHLXMIRU"
$ws.Range("G442").Value = 210.2

$ws.Range("B443").Value = "MICEL"
$ws.Range("G443").Value = 50

$ws.Range("B444").Value = "NDTV"
$ws.Range("G444").Value = 229.35

$ws.Range("B445").Value = "POWERGRID"
$ws.Range("G445").Value = 312.8

$ws.Range("B446").Value = "PRECAM"
$ws.Range("G446").Value = 184.85

$ws.Range("B447").Value = "PRINCEPIPE"
$ws.Range("G447").Value = 614.3

$ws.Range("B448").Value = "RADIANTCMS"
$ws.Range("G448").Value = 79

$ws.Range("B449").Value = "RADICO"
$ws.Range("G449").Value = 1622.6

$ws.Range("B450").Value = "RAMCOCEM"
$ws.Range("G450").Value = 762.85

$ws.Range("B451").Value = "RPSGVENT"
$ws.Range("G451").Value = 641.75

$ws.Range("B452").Value = "RUPA"
$ws.Range("G452").Value = 248.5

$ws.Range("B453").Value = "SALASAR"
$ws.Range("G453").Value = 20.35

$ws.Range("A454").Value = "28/05/2024"
